$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1587, 1).Value = 1586
$ws.Cells.Item(1587, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1587, 3).Value = "5:33 AM"
$ws.Cells.Item(1587, 4).Value = "LO6288"
$ws.Cells.Item(1587, 5).Value = "Zanzibar"
$ws.Cells.Item(1587, 6).Value = "(ZNZ)"
$ws.Cells.Item(1587, 7).Value = "LOT "
$ws.Cells.Item(1587, 8).Value = "B789"
$ws.Cells.Item(1587, 9).Value = "(SP-LSG)"
$ws.Cells.Item(1587, 10).Value = "5:35 AM"
$ws.Cells.Item(1587, 12).Value = "0 hours, 2 minutes"

$ws.Cells.Item(1588, 1).Value = 1587
$ws.Cells.Item(1588, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1588, 3).Value = "5:39 AM"
$ws.Cells.Item(1588, 4).Value = "QY5546"
$ws.Cells.Item(1588, 5).Value = "Leipzig"
$ws.Cells.Item(1588, 6).Value = "(LEJ)"
$ws.Cells.Item(1588, 7).Value = "DHL "
$ws.Cells.Item(1588, 8).Value = "A306"
$ws.Cells.Item(1588, 9).Value = "(D-AEAN)"
$ws.Cells.Item(1588, 10).Value = "5:31 AM"
$ws.Cells.Item(1588, 12).Value = "0 hours, -8 minutes"

$ws.Cells.Item(1589, 1).Value = 1588
$ws.Cells.Item(1589, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1589, 3).Value = "5:49 AM"
$ws.Cells.Item(1589, 4).Value = "5X272"
$ws.Cells.Item(1589, 5).Value = "Cologne"
$ws.Cells.Item(1589, 6).Value = "(CGN)"
$ws.Cells.Item(1589, 7).Value = "UPS "
$ws.Cells.Item(1589, 8).Value = "B752"
$ws.Cells.Item(1589, 9).Value = "(N431UP)"
$ws.Cells.Item(1589, 10).Value = "5:28 AM"
$ws.Cells.Item(1589, 12).Value = "0 hours, -21 minutes"

$ws.Cells.Item(1590, 1).Value = 1589
$ws.Cells.Item(1590, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1590, 3).Value = "6:10 AM"
$ws.Cells.Item(1590, 4).Value = "W61540"
$ws.Cells.Item(1590, 5).Value = "Reykjavik"
$ws.Cells.Item(1590, 6).Value = "(KEF)"
$ws.Cells.Item(1590, 7).Value = "Wizz Air "
$ws.Cells.Item(1590, 8).Value = "A21N"
$ws.Cells.Item(1590, 9).Value = "(HA-LZE)"
$ws.Cells.Item(1590, 10).Value = "6:26 AM"
$ws.Cells.Item(1590, 12).Value = "0 hours, 16 minutes"

$ws.Cells.Item(1591, 1).Value = 1590
$ws.Cells.Item(1591, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1591, 3).Value = "6:15 AM"
$ws.Cells.Item(1591, 4).Value = "LO728"
$ws.Cells.Item(1591, 5).Value = "Yerevan"
$ws.Cells.Item(1591, 6).Value = "(EVN)"
$ws.Cells.Item(1591, 7).Value = "LOT "
$ws.Cells.Item(1591, 8).Value = "B738"
$ws.Cells.Item(1591, 9).Value = "(SP-LWA)"
$ws.Cells.Item(1591, 10).Value = "6:37 AM"
$ws.Cells.Item(1591, 12).Value = "0 hours, 22 minutes"

$ws.Cells.Item(1592, 1).Value = 1591
$ws.Cells.Item(1592, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1592, 3).Value = "6:15 AM"
$ws.Cells.Item(1592, 4).Value = "LO720"
$ws.Cells.Item(1592, 5).Value = "Baku"
$ws.Cells.Item(1592, 6).Value = "(GYD)"
$ws.Cells.Item(1592, 7).Value = "LOT "
$ws.Cells.Item(1592, 8).Value = "B38M"
$ws.Cells.Item(1592, 9).Value = "(SP-LVA)"
$ws.Cells.Item(1592, 10).Value = "6:39 AM"
$ws.Cells.Item(1592, 12).Value = "0 hours, 24 minutes"

$ws.Cells.Item(1593, 1).Value = 1592
$ws.Cells.Item(1593, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1593, 3).Value = "6:25 AM"
$ws.Cells.Item(1593, 4).Value = "LO504"
$ws.Cells.Item(1593, 5).Value = "Ostrava"
$ws.Cells.Item(1593, 6).Value = "(OSR)"
$ws.Cells.Item(1593, 7).Value = "LOT "
$ws.Cells.Item(1593, 8).Value = "E170"
$ws.Cells.Item(1593, 9).Value = "(SP-LDG)"
$ws.Cells.Item(1593, 10).Value = "6:03 AM"
$ws.Cells.Item(1593, 12).Value = "0 hours, -22 minutes"

$ws.Cells.Item(1594, 1).Value = 1593
$ws.Cells.Item(1594, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1594, 3).Value = "6:25 AM"
$ws.Cells.Item(1594, 4).Value = "LO774"
$ws.Cells.Item(1594, 5).Value = "Vilnius"
$ws.Cells.Item(1594, 6).Value = "(VNO)"
$ws.Cells.Item(1594, 7).Value = "LOT "
$ws.Cells.Item(1594, 8).Value = "E195"
$ws.Cells.Item(1594, 9).Value = "(SP-LNO)"
$ws.Cells.Item(1594, 10).Value = "6:17 AM"
$ws.Cells.Item(1594, 12).Value = "0 hours, -8 minutes"

$ws.Cells.Item(1595, 1).Value = 1594
$ws.Cells.Item(1595, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1595, 3).Value = "6:25 AM"
$ws.Cells.Item(1595, 4).Value = "LO644"
$ws.Cells.Item(1595, 5).Value = "Bucharest"
$ws.Cells.Item(1595, 6).Value = "(OTP)"
$ws.Cells.Item(1595, 7).Value = "LOT "
$ws.Cells.Item(1595, 8).Value = "E190"
$ws.Cells.Item(1595, 9).Value = "(SP-LMB)"
$ws.Cells.Item(1595, 10).Value = "6:05 AM"
$ws.Cells.Item(1595, 12).Value = "0 hours, -20 minutes"

$ws.Cells.Item(1596, 1).Value = 1595
$ws.Cells.Item(1596, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1596, 3).Value = "6:25 AM"
$ws.Cells.Item(1596, 4).Value = "LO792"
$ws.Cells.Item(1596, 5).Value = "Tallinn"
$ws.Cells.Item(1596, 6).Value = "(TLL)"
$ws.Cells.Item(1596, 7).Value = "LOT "
$ws.Cells.Item(1596, 8).Value = "E190"
$ws.Cells.Item(1596, 9).Value = "(SP-LMA)"
$ws.Cells.Item(1596, 10).Value = "6:22 AM"
$ws.Cells.Item(1596, 12).Value = "0 hours, -3 minutes"

$ws.Cells.Item(1597, 1).Value = 1596
$ws.Cells.Item(1597, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1597, 3).Value = "6:30 AM"
$ws.Cells.Item(1597, 4).Value = "LO794"
$ws.Cells.Item(1597, 5).Value = "Riga"
$ws.Cells.Item(1597, 6).Value = "(RIX)"
$ws.Cells.Item(1597, 7).Value = "LOT "
$ws.Cells.Item(1597, 8).Value = "E75S"
$ws.Cells.Item(1597, 9).Value = "(SP-LIN)"
$ws.Cells.Item(1597, 10).Value = "6:20 AM"
$ws.Cells.Item(1597, 12).Value = "0 hours, -10 minutes"

$ws.Cells.Item(1598, 1).Value = 1597
$ws.Cells.Item(1598, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1598, 3).Value = "6:35 AM"
$ws.Cells.Item(1598, 4).Value = "LO3828"
$ws.Cells.Item(1598, 5).Value = "Gdansk"
$ws.Cells.Item(1598, 6).Value = "(GDN)"
$ws.Cells.Item(1598, 7).Value = "LOT "
$ws.Cells.Item(1598, 8).Value = "E195"
$ws.Cells.Item(1598, 9).Value = "(SP-LNK)"
$ws.Cells.Item(1598, 10).Value = "6:29 AM"
$ws.Cells.Item(1598, 12).Value = "0 hours, -6 minutes"

$ws.Cells.Item(1599, 1).Value = 1598
$ws.Cells.Item(1599, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1599, 3).Value = "6:35 AM"
$ws.Cells.Item(1599, 4).Value = "LO3850"
$ws.Cells.Item(1599, 5).Value = "Wroclaw"
$ws.Cells.Item(1599, 6).Value = "(WRO)"
$ws.Cells.Item(1599, 7).Value = "LOT "
$ws.Cells.Item(1599, 8).Value = "E75S"
$ws.Cells.Item(1599, 9).Value = "(SP-LIK)"
$ws.Cells.Item(1599, 10).Value = "6:13 AM"
$ws.Cells.Item(1599, 12).Value = "0 hours, -22 minutes"

$ws.Cells.Item(1600, 1).Value = 1599
$ws.Cells.Item(1600, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1600, 3).Value = "6:35 AM"
$ws.Cells.Item(1600, 4).Value = "LO3880"
$ws.Cells.Item(1600, 5).Value = "Katowice"
$ws.Cells.Item(1600, 6).Value = "(KTW)"
$ws.Cells.Item(1600, 7).Value = "LOT "
$ws.Cells.Item(1600, 8).Value = "E190"
$ws.Cells.Item(1600, 9).Value = "(SP-LMG)"
$ws.Cells.Item(1600, 10).Value = "6:33 AM"
$ws.Cells.Item(1600, 12).Value = "0 hours, -2 minutes"

$ws.Cells.Item(1601, 1).Value = 1600
$ws.Cells.Item(1601, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1601, 3).Value = "6:35 AM"
$ws.Cells.Item(1601, 4).Value = "LO3910"
$ws.Cells.Item(1601, 5).Value = "Krakow"
$ws.Cells.Item(1601, 6).Value = "(KRK)"
$ws.Cells.Item(1601, 7).Value = "LOT "
$ws.Cells.Item(1601, 8).Value = "E195"
$ws.Cells.Item(1601, 9).Value = "(SP-LNM)"
$ws.Cells.Item(1601, 10).Value = "6:30 AM"
$ws.Cells.Item(1601, 12).Value = "0 hours, -5 minutes"

$ws.Cells.Item(1602, 1).Value = 1601
$ws.Cells.Item(1602, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1602, 3).Value = "6:40 AM"
$ws.Cells.Item(1602, 4).Value = "LO3804"
$ws.Cells.Item(1602, 5).Value = "Rzeszow"
$ws.Cells.Item(1602, 6).Value = "(RZE)"
$ws.Cells.Item(1602, 7).Value = "LOT (Sliwka Naleczowska Livery) "
$ws.Cells.Item(1602, 8).Value = "E195"
$ws.Cells.Item(1602, 9).Value = "(SP-LNC)"
$ws.Cells.Item(1602, 10).Value = "6:15 AM"
$ws.Cells.Item(1602, 12).Value = "0 hours, -25 minutes"

$ws.Cells.Item(1603, 1).Value = 1602
$ws.Cells.Item(1603, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1603, 3).Value = "6:40 AM"
$ws.Cells.Item(1603, 4).Value = "LO3942"
$ws.Cells.Item(1603, 5).Value = "Poznan"
$ws.Cells.Item(1603, 6).Value = "(POZ)"
$ws.Cells.Item(1603, 7).Value = "LOT "
$ws.Cells.Item(1603, 8).Value = "E190"
$ws.Cells.Item(1603, 9).Value = "(SP-LMH)"
$ws.Cells.Item(1603, 10).Value = "6:35 AM"
$ws.Cells.Item(1603, 12).Value = "0 hours, -5 minutes"

$ws.Cells.Item(1604, 1).Value = 1603
$ws.Cells.Item(1604, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1604, 3).Value = "7:25 AM"
$ws.Cells.Item(1604, 4).Value = "LO726"
$ws.Cells.Item(1604, 5).Value = "Tbilisi"
$ws.Cells.Item(1604, 6).Value = "(TBS)"
$ws.Cells.Item(1604, 7).Value = "LOT "
$ws.Cells.Item(1604, 8).Value = "B38M"
$ws.Cells.Item(1604, 9).Value = "(SP-LVB)"
$ws.Cells.Item(1604, 10).Value = "7:19 AM"
$ws.Cells.Item(1604, 12).Value = "0 hours, -6 minutes"

$ws.Cells.Item(1605, 1).Value = 1604
$ws.Cells.Item(1605, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1605, 3).Value = "8:05 AM"
$ws.Cells.Item(1605, 4).Value = "TK4083"
$ws.Cells.Item(1605, 5).Value = "Antalya"
$ws.Cells.Item(1605, 6).Value = "(AYT)"
$ws.Cells.Item(1605, 7).Value = "Turkish Airlines "
$ws.Cells.Item(1605, 8).Value = "B738"
$ws.Cells.Item(1605, 9).Value = "(TC-JVO)"
$ws.Cells.Item(1605, 10).Value = "8:03 AM"
$ws.Cells.Item(1605, 12).Value = "0 hours, -2 minutes"

$ws.Cells.Item(1606, 1).Value = 1605
$ws.Cells.Item(1606, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1606, 3).Value = "8:15 AM"
$ws.Cells.Item(1606, 4).Value = "AY1141"
$ws.Cells.Item(1606, 5).Value = "Helsinki"
$ws.Cells.Item(1606, 6).Value = "(HEL)"
$ws.Cells.Item(1606, 7).Value = "Finnair (Oneworld livery) "
$ws.Cells.Item(1606, 8).Value = "E190"
$ws.Cells.Item(1606, 9).Value = "(OH-LKN)"
$ws.Cells.Item(1606, 10).Value = "8:27 AM"
$ws.Cells.Item(1606, 12).Value = "0 hours, 12 minutes"

$ws.Cells.Item(1607, 1).Value = 1606
$ws.Cells.Item(1607, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1607, 3).Value = "8:25 AM"
$ws.Cells.Item(1607, 4).Value = "LO152"
$ws.Cells.Item(1607, 5).Value = "Tel Aviv"
$ws.Cells.Item(1607, 6).Value = "(TLV)"
$ws.Cells.Item(1607, 7).Value = "LOT "
$ws.Cells.Item(1607, 8).Value = "B38M"
$ws.Cells.Item(1607, 9).Value = "(SP-LVF)"
$ws.Cells.Item(1607, 10).Value = "8:14 AM"
$ws.Cells.Item(1607, 12).Value = "0 hours, -11 minutes"

$ws.Cells.Item(1608, 1).Value = 1607
$ws.Cells.Item(1608, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1608, 3).Value = "8:25 AM"
$ws.Cells.Item(1608, 4).Value = "LO402"
$ws.Cells.Item(1608, 5).Value = "Hamburg"
$ws.Cells.Item(1608, 6).Value = "(HAM)"
$ws.Cells.Item(1608, 7).Value = "LOT "
$ws.Cells.Item(1608, 8).Value = "E170"
$ws.Cells.Item(1608, 9).Value = "(SP-LDI)"
$ws.Cells.Item(1608, 10).Value = "8:13 AM"
$ws.Cells.Item(1608, 12).Value = "0 hours, -12 minutes"

$ws.Cells.Item(1609, 1).Value = 1608
$ws.Cells.Item(1609, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1609, 3).Value = "8:25 AM"
$ws.Cells.Item(1609, 4).Value = "LO460"
$ws.Cells.Item(1609, 5).Value = "Copenhagen"
$ws.Cells.Item(1609, 6).Value = "(CPH)"
$ws.Cells.Item(1609, 7).Value = "LOT "
$ws.Cells.Item(1609, 8).Value = "E195"
$ws.Cells.Item(1609, 9).Value = "(SP-LNN)"
$ws.Cells.Item(1609, 10).Value = "8:08 AM"
$ws.Cells.Item(1609, 12).Value = "0 hours, -17 minutes"

$ws.Cells.Item(1610, 1).Value = 1609
$ws.Cells.Item(1610, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1610, 3).Value = "8:35 AM"
$ws.Cells.Item(1610, 4).Value = "LO138"
$ws.Cells.Item(1610, 5).Value = "Istanbul"
$ws.Cells.Item(1610, 6).Value = "(IST)"
$ws.Cells.Item(1610, 7).Value = "LOT "
$ws.Cells.Item(1610, 8).Value = "E195"
$ws.Cells.Item(1610, 9).Value = "(SP-LNL)"
$ws.Cells.Item(1610, 10).Value = "8:21 AM"
$ws.Cells.Item(1610, 12).Value = "0 hours, -14 minutes"

$ws.Cells.Item(1611, 1).Value = 1610
$ws.Cells.Item(1611, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1611, 3).Value = "8:35 AM"
$ws.Cells.Item(1611, 4).Value = "LO530"
$ws.Cells.Item(1611, 5).Value = "Prague"
$ws.Cells.Item(1611, 6).Value = "(PRG)"
$ws.Cells.Item(1611, 7).Value = "LOT "
$ws.Cells.Item(1611, 8).Value = "E170"
$ws.Cells.Item(1611, 9).Value = "(SP-LDF)"
$ws.Cells.Item(1611, 10).Value = "8:11 AM"
$ws.Cells.Item(1611, 12).Value = "0 hours, -24 minutes"

$ws.Cells.Item(1612, 1).Value = 1611
$ws.Cells.Item(1612, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1612, 3).Value = "8:35 AM"
$ws.Cells.Item(1612, 4).Value = "TK1265"
$ws.Cells.Item(1612, 5).Value = "Istanbul"
$ws.Cells.Item(1612, 6).Value = "(IST)"
$ws.Cells.Item(1612, 7).Value = "Turkish Airlines "
$ws.Cells.Item(1612, 8).Value = "A321"
$ws.Cells.Item(1612, 9).Value = "(TC-JSE)"
$ws.Cells.Item(1612, 10).Value = "8:25 AM"
$ws.Cells.Item(1612, 12).Value = "0 hours, -10 minutes"

$ws.Cells.Item(1613, 1).Value = 1612
$ws.Cells.Item(1613, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1613, 3).Value = "8:40 AM"
$ws.Cells.Item(1613, 4).Value = "LH1346"
$ws.Cells.Item(1613, 5).Value = "Frankfurt"
$ws.Cells.Item(1613, 6).Value = "(FRA)"
$ws.Cells.Item(1613, 7).Value = "Lufthansa "
$ws.Cells.Item(1613, 8).Value = "A320"
$ws.Cells.Item(1613, 9).Value = "(D-AIUK)"
$ws.Cells.Item(1613, 10).Value = "8:38 AM"
$ws.Cells.Item(1613, 12).Value = "0 hours, -2 minutes"

$ws.Cells.Item(1614, 1).Value = 1613
$ws.Cells.Item(1614, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1614, 3).Value = "9:00 AM"
$ws.Cells.Item(1614, 4).Value = "LX1342"
$ws.Cells.Item(1614, 5).Value = "Zurich"
$ws.Cells.Item(1614, 6).Value = "(ZRH)"
$ws.Cells.Item(1614, 7).Value = "Helvetic Airways "
$ws.Cells.Item(1614, 8).Value = "E290"
$ws.Cells.Item(1614, 9).Value = "(HB-AZB)"
$ws.Cells.Item(1614, 10).Value = "8:46 AM"
$ws.Cells.Item(1614, 12).Value = "0 hours, -14 minutes"

$ws.Cells.Item(1615, 1).Value = 1614
$ws.Cells.Item(1615, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1615, 3).Value = "9:05 AM"
$ws.Cells.Item(1615, 4).Value = "LO780"
$ws.Cells.Item(1615, 5).Value = "Vilnius"
$ws.Cells.Item(1615, 6).Value = "(VNO)"
$ws.Cells.Item(1615, 7).Value = "LOT "
$ws.Cells.Item(1615, 8).Value = "E190"
$ws.Cells.Item(1615, 9).Value = "(SP-LMD)"
$ws.Cells.Item(1615, 10).Value = "8:58 AM"
$ws.Cells.Item(1615, 12).Value = "0 hours, -7 minutes"

$ws.Cells.Item(1616, 1).Value = 1615
$ws.Cells.Item(1616, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1616, 3).Value = "9:10 AM"
$ws.Cells.Item(1616, 4).Value = "LO384"
$ws.Cells.Item(1616, 5).Value = "Frankfurt"
$ws.Cells.Item(1616, 6).Value = "(FRA)"
$ws.Cells.Item(1616, 7).Value = "LOT "
$ws.Cells.Item(1616, 8).Value = "E75S"
$ws.Cells.Item(1616, 9).Value = "(SP-LIC)"
$ws.Cells.Item(1616, 10).Value = "8:44 AM"
$ws.Cells.Item(1616, 12).Value = "0 hours, -26 minutes"

$ws.Cells.Item(1617, 1).Value = 1616
$ws.Cells.Item(1617, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1617, 3).Value = "9:10 AM"
$ws.Cells.Item(1617, 4).Value = "LO406"
$ws.Cells.Item(1617, 5).Value = "Dusseldorf"
$ws.Cells.Item(1617, 6).Value = "(DUS)"
$ws.Cells.Item(1617, 7).Value = "LOT "
$ws.Cells.Item(1617, 8).Value = "E75S"
$ws.Cells.Item(1617, 9).Value = "(SP-LID)"
$ws.Cells.Item(1617, 10).Value = "8:52 AM"
$ws.Cells.Item(1617, 12).Value = "0 hours, -18 minutes"

$ws.Cells.Item(1618, 1).Value = 1617
$ws.Cells.Item(1618, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1618, 3).Value = "9:15 AM"
$ws.Cells.Item(1618, 4).Value = "SK2601"
$ws.Cells.Item(1618, 5).Value = "Stockholm"
$ws.Cells.Item(1618, 6).Value = "(ARN)"
$ws.Cells.Item(1618, 7).Value = "SAS "
$ws.Cells.Item(1618, 8).Value = "CRJ9"
$ws.Cells.Item(1618, 9).Value = "(ES-ACJ)"
$ws.Cells.Item(1618, 10).Value = "9:34 AM"
$ws.Cells.Item(1618, 12).Value = "0 hours, 19 minutes"

$ws.Cells.Item(1619, 1).Value = 1618
$ws.Cells.Item(1619, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1619, 3).Value = "9:25 AM"
$ws.Cells.Item(1619, 4).Value = "SK751"
$ws.Cells.Item(1619, 5).Value = "Copenhagen"
$ws.Cells.Item(1619, 6).Value = "(CPH)"
$ws.Cells.Item(1619, 7).Value = "SAS "
$ws.Cells.Item(1619, 8).Value = "CRJ9"
$ws.Cells.Item(1619, 9).Value = "(ES-ACK)"
$ws.Cells.Item(1619, 10).Value = "9:29 AM"
$ws.Cells.Item(1619, 12).Value = "0 hours, 4 minutes"

$ws.Cells.Item(1620, 1).Value = 1619
$ws.Cells.Item(1620, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1620, 3).Value = "9:30 AM"
$ws.Cells.Item(1620, 4).Value = "LO456"
$ws.Cells.Item(1620, 5).Value = "Stockholm"
$ws.Cells.Item(1620, 6).Value = "(ARN)"
$ws.Cells.Item(1620, 7).Value = "LOT (Grzeski Livery) "
$ws.Cells.Item(1620, 8).Value = "E195"
$ws.Cells.Item(1620, 9).Value = "(SP-LNB)"
$ws.Cells.Item(1620, 10).Value = "9:18 AM"
$ws.Cells.Item(1620, 12).Value = "0 hours, -12 minutes"

$ws.Cells.Item(1621, 1).Value = 1620
$ws.Cells.Item(1621, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1621, 3).Value = "9:30 AM"
$ws.Cells.Item(1621, 4).Value = "OS625"
$ws.Cells.Item(1621, 5).Value = "Vienna"
$ws.Cells.Item(1621, 6).Value = "(VIE)"
$ws.Cells.Item(1621, 7).Value = "Austrian Airlines "
$ws.Cells.Item(1621, 8).Value = "E195"
$ws.Cells.Item(1621, 9).Value = "(OE-LWO)"
$ws.Cells.Item(1621, 10).Value = "9:10 AM"
$ws.Cells.Item(1621, 12).Value = "0 hours, -20 minutes"

$ws.Cells.Item(1622, 1).Value = 1621
$ws.Cells.Item(1622, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1622, 3).Value = "9:35 AM"
$ws.Cells.Item(1622, 4).Value = "LO496"
$ws.Cells.Item(1622, 5).Value = "Gothenburg"
$ws.Cells.Item(1622, 6).Value = "(GOT)"
$ws.Cells.Item(1622, 7).Value = "LOT "
$ws.Cells.Item(1622, 8).Value = "E195"
$ws.Cells.Item(1622, 9).Value = "(SP-LNG)"
$ws.Cells.Item(1622, 10).Value = "9:21 AM"
$ws.Cells.Item(1622, 12).Value = "0 hours, -14 minutes"

$ws.Cells.Item(1623, 1).Value = 1622
$ws.Cells.Item(1623, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1623, 3).Value = "9:40 AM"
$ws.Cells.Item(1623, 4).Value = "LO232"
$ws.Cells.Item(1623, 5).Value = "Brussels"
$ws.Cells.Item(1623, 6).Value = "(BRU)"
$ws.Cells.Item(1623, 7).Value = "LOT (Retro Livery) "
$ws.Cells.Item(1623, 8).Value = "E75S"
$ws.Cells.Item(1623, 9).Value = "(SP-LIM)"
$ws.Cells.Item(1623, 10).Value = "9:13 AM"
$ws.Cells.Item(1623, 12).Value = "0 hours, -27 minutes"

$ws.Cells.Item(1624, 1).Value = 1623
$ws.Cells.Item(1624, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1624, 3).Value = "9:45 AM"
$ws.Cells.Item(1624, 4).Value = "LO3852"
$ws.Cells.Item(1624, 5).Value = "Wroclaw"
$ws.Cells.Item(1624, 6).Value = "(WRO)"
$ws.Cells.Item(1624, 7).Value = "LOT "
$ws.Cells.Item(1624, 8).Value = "E75S"
$ws.Cells.Item(1624, 9).Value = "(SP-LIB)"
$ws.Cells.Item(1624, 10).Value = "9:23 AM"
$ws.Cells.Item(1624, 12).Value = "0 hours, -22 minutes"

$ws.Cells.Item(1625, 1).Value = 1624
$ws.Cells.Item(1625, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1625, 3).Value = "9:45 AM"
$ws.Cells.Item(1625, 4).Value = "LO3904"
$ws.Cells.Item(1625, 5).Value = "Krakow"
$ws.Cells.Item(1625, 6).Value = "(KRK)"
$ws.Cells.Item(1625, 7).Value = "LOT "
$ws.Cells.Item(1625, 8).Value = "E75S"
$ws.Cells.Item(1625, 9).Value = "(SP-LIA)"
$ws.Cells.Item(1625, 10).Value = "9:31 AM"
$ws.Cells.Item(1625, 12).Value = "0 hours, -14 minutes"

$ws.Cells.Item(1626, 1).Value = 1625
$ws.Cells.Item(1626, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1626, 3).Value = "9:45 AM"
$ws.Cells.Item(1626, 4).Value = "LO3982"
$ws.Cells.Item(1626, 5).Value = "Zielona Gora"
$ws.Cells.Item(1626, 6).Value = "(IEG)"
$ws.Cells.Item(1626, 7).Value = "LOT "
$ws.Cells.Item(1626, 8).Value = "E75S"
$ws.Cells.Item(1626, 9).Value = "(SP-LIL)"
$ws.Cells.Item(1626, 10).Value = "9:39 AM"
$ws.Cells.Item(1626, 12).Value = "0 hours, -6 minutes"

$ws.Cells.Item(1627, 1).Value = 1626
$ws.Cells.Item(1627, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1627, 3).Value = "9:50 AM"
$ws.Cells.Item(1627, 4).Value = "LO3832"
$ws.Cells.Item(1627, 5).Value = "Gdansk"
$ws.Cells.Item(1627, 6).Value = "(GDN)"
$ws.Cells.Item(1627, 7).Value = "LOT "
$ws.Cells.Item(1627, 8).Value = "E190"
$ws.Cells.Item(1627, 9).Value = "(SP-LME)"
$ws.Cells.Item(1627, 10).Value = "9:44 AM"
$ws.Cells.Item(1627, 12).Value = "0 hours, -6 minutes"

$ws.Cells.Item(1628, 1).Value = 1627
$ws.Cells.Item(1628, 2).Value = "Monday, Jan 16"
$ws.Cells.Item(1628, 3).Value = "10:05 AM"
$ws.Cells.Item(1628, 4).Value = "W61502"
$ws.Cells.Item(1628, 5).Value = "Stockholm"
$ws.Cells.Item(1628, 6).Value = "(NYO)"
$ws.Cells.Item(1628, 7).Value = "Wizz Air "
$ws.Cells.Item(1628, 8).Value = "A321"
$ws.Cells.Item(1628, 9).Value = "(HA-LTB)"
$ws.Cells.Item(1628, 10).Value = "9:42 AM"
$ws.Cells.Item(1628, 12).Value = "0 hours, -23 minutes"
